# Append two new ingredient rows (haselnusskerne, mandelnkerne) to the keto sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "fat/carbs/protein" columns in this sheet store numeric-looking values
# as plain text (shared strings), matching the rest of the existing data.
# Force text number format so Excel doesn't auto-convert them to numbers.
$ws.Range("C51:E52").NumberFormat = "@"

# Row 51: haselnusskerne
$ws.Range("A51").Value = "haselnusskerne"
$ws.Range("B51").Value = 100
$ws.Range("C51").Value = "68.4"
$ws.Range("D51").Value = "5.6"
$ws.Range("E51").Value = "16.3"
$ws.Range("F51").Value = 716

# Row 52: mandelnkerne
$ws.Range("A52").Value = "mandelnkerne"
$ws.Range("B52").Value = 100
$ws.Range("C52").Value = "53.3"
$ws.Range("D52").Value = "4.8"
$ws.Range("E52").Value = "24.5"
$ws.Range("F52").Value = 621

# Revert to the sheet's normal (default) cell style now that the text has
# been entered, so the new cells don't carry an explicit text-format style
# (matches how the rest of the sheet's cells are styled).
$ws.Range("C51:E52").Style = "Standard"

# Update selection to match the author's saved view state (Excel records the
# last selected cell on save).
$ws.Range("F53").Select() | Out-Null
